# Apply crypto price/volume updates per commit "Updated cryptos list on Fri Sep 20 13:23:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.985.95'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.536.04'
$ws.Range('E3').Value = '  +3.97%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '569.75'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.65'
$ws.Range('E6').Value = '  +5.42%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.582'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('D9').Value = '2.534.74'
$ws.Range('E9').Value = '  +3.96%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.105'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.65'
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.355'
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.86'
$ws.Range('E14').Value = '  +5.65%  '
$ws.Range('D15').Value = '2.989.51'
$ws.Range('E15').Value = '  +4.02%  '
$ws.Range('D16').Value = '62.971.65'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000142'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '2.535.42'
$ws.Range('E18').Value = '  +3.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.55'
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '337.19'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.74'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.69'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  -3.40%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.59'
$ws.Range('E26').Value = '  +2.95%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.52'
$ws.Range('E27').Value = '  +11.04%  '
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.37'
$ws.Range('E29').Value = '  +1.51%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.26'
$ws.Range('E30').Value = '  +10.85%  '
$ws.Range('D31').Value = '0.0₃0812'
$ws.Range('E31').Value = '  +2.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.84'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '177.79'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  +7.13%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '416.11'
$ws.Range('E35').Value = '  +9.39%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.400'
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.83'
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.40'
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').Value = '  +1.26%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.32'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '152.54'
$ws.Range('E43').Value = '  +4.89%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.75'
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.70'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.605'
$ws.Range('E46').Value = '  +1.67%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0962'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0520'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0237'
$ws.Range('E49').Value = '  +6.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.41'
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.78'
$ws.Range('E51').Value = '  +1.69%  '
